$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the default (unstyled) cell
# format so numeric-looking strings (e.g. "607.43") are not coerced into
# floating point numbers by Excel's automatic type inference.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '67.001.55'
Set-TextValue $ws.Range("E2") '  +3.18%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.216.33'
Set-TextValue $ws.Range("E3") '  +2.04%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.00%  '

# Row 5
Set-TextValue $ws.Range("D5") '607.43'
Set-TextValue $ws.Range("E5") '  +4.27%  '

# Row 6
Set-TextValue $ws.Range("D6") '158.74'
Set-TextValue $ws.Range("E6") '  +6.50%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.999'
Set-TextValue $ws.Range("E7") '  -0.05%  '

# Row 8
Set-TextValue $ws.Range("B8") 'LidoStakedEther'
Set-TextValue $ws.Range("C8") 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextValue $ws.Range("D8") '3.216.61'
Set-TextValue $ws.Range("E8") '  +2.06%  '

# Row 9
Set-TextValue $ws.Range("B9") 'XRP'
Set-TextValue $ws.Range("C9") 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range("D9") '0.556'
Set-TextValue $ws.Range("E9") '  +5.93%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.162'
Set-TextValue $ws.Range("E10") '  +1.89%  '

# Row 11
Set-TextValue $ws.Range("D11") '6.02'
Set-TextValue $ws.Range("E11") '  -2.88%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.518'
Set-TextValue $ws.Range("E12") '  +3.68%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +1.33%  '

# Row 14
Set-TextValue $ws.Range("D14") '39.67'
Set-TextValue $ws.Range("E14") '  +6.91%  '

# Row 15
Set-TextValue $ws.Range("D15") '3.741.42'
Set-TextValue $ws.Range("E15") '  +2.00%  '

# Row 16
Set-TextValue $ws.Range("D16") '66.992.18'
Set-TextValue $ws.Range("E16") '  +3.23%  '

# Row 17
Set-TextValue $ws.Range("D17") '7.50'
Set-TextValue $ws.Range("E17") '  +5.11%  '

# Row 18
Set-TextValue $ws.Range("D18") '3.213.64'
Set-TextValue $ws.Range("E18") '  +1.94%  '

# Row 19
Set-TextValue $ws.Range("E19") '  +0.95%  '

# Row 20
Set-TextValue $ws.Range("D20") '521.68'
Set-TextValue $ws.Range("E20") '  +3.42%  '

# Row 21
Set-TextValue $ws.Range("D21") '15.48'
Set-TextValue $ws.Range("E21") '  +2.70%  '

# Row 22
Set-TextValue $ws.Range("D22") '0.747'
Set-TextValue $ws.Range("E22") '  +4.71%  '

# Row 23
Set-TextValue $ws.Range("D23") '8.25'
Set-TextValue $ws.Range("E23") '  +6.46%  '

# Row 24
Set-TextValue $ws.Range("D24") '15.22'
Set-TextValue $ws.Range("E24") '  +1.08%  '

# Row 25
Set-TextValue $ws.Range("D25") '85.58'
Set-TextValue $ws.Range("E25") '  +1.56%  '

# Row 26
Set-TextValue $ws.Range("D26") '1.00'
Set-TextValue $ws.Range("E26") '  +0.01%  '

# Row 27
Set-TextValue $ws.Range("D27") '9.48'
Set-TextValue $ws.Range("E27") '  +4.92%  '

# Row 28
Set-TextValue $ws.Range("D28") '3.06'
Set-TextValue $ws.Range("E28") '  +4.80%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.43'
Set-TextValue $ws.Range("E29") '  +11.41%  '

# Row 30
Set-TextValue $ws.Range("D30") '3.07'
Set-TextValue $ws.Range("E30") '  +9.50%  '

# Row 31
Set-TextValue $ws.Range("D31") '7.02'
Set-TextValue $ws.Range("E31") '  +9.04%  '

# Row 32
Set-TextValue $ws.Range("D32") '28.48'

# Row 33
Set-TextValue $ws.Range("D33") '1.25'
Set-TextValue $ws.Range("E33") '  +3.46%  '

# Row 34
Set-TextValue $ws.Range("E34") '  +0.16%  '

# Row 35
Set-TextValue $ws.Range("D35") '6.64'
Set-TextValue $ws.Range("E35") '  +2.38%  '

# Row 36
Set-TextValue $ws.Range("D36") '527.09'
Set-TextValue $ws.Range("E36") '  +11.13%  '

# Row 37
Set-TextValue $ws.Range("D37") '55.18'
Set-TextValue $ws.Range("E37") '  +0.73%  '

# Row 38
Set-TextValue $ws.Range("D38") '0.0913'
Set-TextValue $ws.Range("E38") '  +2.13%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.0430'
Set-TextValue $ws.Range("E39") '  +3.96%  '

# Row 40
Set-TextValue $ws.Range("D40") '0.128'
Set-TextValue $ws.Range("E40") '  +9.64%  '

# Row 41
Set-TextValue $ws.Range("D41") '8.98'
Set-TextValue $ws.Range("E41") '  +3.18%  '

# Row 42
Set-TextValue $ws.Range("D42") '2.95'
Set-TextValue $ws.Range("E42") '  +0.51%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.0₃0697'
Set-TextValue $ws.Range("E43") '  +15.98%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.306'
Set-TextValue $ws.Range("E44") '  +8.90%  '

# Row 45
Set-TextValue $ws.Range("D45") '2.55'
Set-TextValue $ws.Range("E45") '  +4.92%  '

# Row 46
Set-TextValue $ws.Range("D46") '2.915.74'
Set-TextValue $ws.Range("E46") '  -2.68%  '

# Row 47
Set-TextValue $ws.Range("D47") '29.11'
Set-TextValue $ws.Range("E47") '  +2.54%  '

# Row 48
Set-TextValue $ws.Range("B48") 'ThetaToken'
Set-TextValue $ws.Range("C48") 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range("D48") '2.44'
Set-TextValue $ws.Range("E48") '  +9.12%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.119'
Set-TextValue $ws.Range("E49") '  +4.01%  '

# Row 50
Set-TextValue $ws.Range("B50") 'CoreDAO'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-TextValue $ws.Range("D50") '2.70'
Set-TextValue $ws.Range("E50") '  +9.27%  '

# Row 51
Set-TextValue $ws.Range("B51") 'USDe'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D51") '0.999'
Set-TextValue $ws.Range("E51") '  +0.00%  '
